$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Jun Yang
$ws.Range("B2").Value = "jun@mail.com"
$ws.Range("C2").Value = "junyang"
$ws.Range("D2").Value = "password"
$ws.Range("E2").Value = $false

# Update row 3: Quan
$ws.Range("B3").Value = "quan@mail.com"
$ws.Range("C3").Value = "quan"
$ws.Range("D3").Value = "password"
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = "c9d8e441332d46bbb9655b8239c26e94"

# Remove rows 4 and 5 (James Doe, Ben Doe)
$ws.Range("A4:F5").Delete()
